# "Generate Report for Handoff"
#
# The localization-status report is regenerated: the ff42255e-...-md file
# has now been handed off for localization, so its status moves from
# "Handed back: in sync with en-US" to "Ready for handoff" on every sheet,
# timestamps are refreshed, and the zh-cn / de-de detail sheets record a
# warning that the previous handback file is out of date.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-22 00:58:42"

# --- zh-cn detail sheet -----------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("H3").Value = "2016-08-22 00:58:38"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b7eae0b427bae54264f84c152b3a0a9c22833e8c/e2e/ff42255e-aacc-4a51-bd5d-b20f3d687c0b.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/825724a915c6c63c814c9051e3049ad6630d3263/e2e/ff42255e-aacc-4a51-bd5d-b20f3d687c0b.md."
# widen the Error Detail column so the new message is readable
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664

# --- de-de detail sheet -----------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("H3").Value = "2016-08-22 00:58:42"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b7eae0b427bae54264f84c152b3a0a9c22833e8c/e2e/ff42255e-aacc-4a51-bd5d-b20f3d687c0b.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/825724a915c6c63c814c9051e3049ad6630d3263/e2e/ff42255e-aacc-4a51-bd5d-b20f3d687c0b.md."
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664
